$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Turn the Census 2000 PDF URL (Median Age section) into a real
#    hyperlink, styled with the "Hyperlink" character style.
# ------------------------------------------------------------------
$oldUrl = "https://www2.census.gov/programs-surveys/decennial/2000/phc/phc-t-09/tab07.pdf"
$newUrl = "https://www.census.gov/data/tables/time-series/demo/popest/2020s-national-detail.html"

# Find the paragraph index that currently holds the old URL, before
# any mutation happens (paragraph indices stay valid afterwards).
$paraIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $probe = $d.Range($p.Range.Start, $p.Range.End)
    if ($probe.Text -like "*$oldUrl*") {
        $paraIndex = $i
    }
}

$findRange = $d.Content
$find = $findRange.Find
$find.ClearFormatting()
$found = $find.Execute($oldUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($findRange, $oldUrl) | Out-Null
}

if ($paraIndex -gt 0) {
    # ------------------------------------------------------------------
    # 2. Add a new paragraph right after that one with the new
    #    (2020s) Census source link as plain text.
    # ------------------------------------------------------------------
    $hostPara = $d.Paragraphs.Item($paraIndex)
    $hostPara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $insertPoint.InsertAfter("US Census: ")

    $newPara2 = $d.Paragraphs.Item($paraIndex + 1)
    $urlStart = $newPara2.Range.Start + 11
    $urlPoint = $d.Range($urlStart, $urlStart)
    $urlPoint.InsertAfter($newUrl)
}

# ------------------------------------------------------------------
# 3. Register the "Hyperlink" and "Unresolved Mention" character
#    styles (brought in by Word when the hyperlink was created).
# ------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = "DefaultParagraphFont"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.Font.Color = 0xC16305
$hyperlinkStyle.Font.Underline = 1

$mentionStyle = $d.Styles.Add("Unresolved Mention", 2)
$mentionStyle.BaseStyle = "DefaultParagraphFont"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 0x5C5E60
